$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Capacitor"
$ws.Cells.Item(1,2).Value = "Voltage"
$ws.Cells.Item(1,3).Value = "Temperature"

# ---- Column A (Capacitor) values, first-use order e100v then e300v ----
$ws.Cells.Item(2,1).Value  = "e100v"
$ws.Cells.Item(3,1).Value  = "e300v"
$ws.Cells.Item(4,1).Value  = "e100v"
$ws.Cells.Item(5,1).Value  = "e300v"
$ws.Cells.Item(6,1).Value  = "e100v"
$ws.Cells.Item(7,1).Value  = "e300v"
$ws.Cells.Item(8,1).Value  = "e100v"
$ws.Cells.Item(9,1).Value  = "e300v"
$ws.Cells.Item(10,1).Value = "e100v"
$ws.Cells.Item(11,1).Value = "e300v"
$ws.Cells.Item(12,1).Value = "e100v"
$ws.Cells.Item(13,1).Value = "e300v"

# ---- Header D1 (leakage) ----
$ws.Cells.Item(1,4).Value = "leakage"

# ---- Column D (leakage) values entered in the same order the lab data was recorded ----
$ws.Cells.Item(4,4).Value  = "5uA"
$ws.Cells.Item(2,4).Value  = "12.29uA"
$ws.Cells.Item(3,4).Value  = "7.8uA"
$ws.Cells.Item(5,4).Value  = "3.08uA"
$ws.Cells.Item(9,4).Value  = "12.25uA"
$ws.Cells.Item(8,4).Value  = "4uA"
$ws.Cells.Item(6,4).Value  = "3.093uA"
$ws.Cells.Item(7,4).Value  = "7.6uA"
$ws.Cells.Item(11,4).Value = "10.3uA"
$ws.Cells.Item(10,4).Value = "3.72uA"
$ws.Cells.Item(12,4).Value = "4.3uA"
$ws.Cells.Item(13,4).Value = "17.23uA"

# ---- Column B (Voltage) ----
$ws.Cells.Item(2,2).Value  = 50
$ws.Cells.Item(3,2).Value  = 150
$ws.Cells.Item(4,2).Value  = 62
$ws.Cells.Item(5,2).Value  = 200
$ws.Cells.Item(6,2).Value  = 50
$ws.Cells.Item(7,2).Value  = 150
$ws.Cells.Item(8,2).Value  = 62
$ws.Cells.Item(9,2).Value  = 200
$ws.Cells.Item(10,2).Value = 50
$ws.Cells.Item(11,2).Value = 150
$ws.Cells.Item(12,2).Value = 62
$ws.Cells.Item(13,2).Value = 200

# ---- Column C (Temperature) ----
$ws.Cells.Item(2,3).Value  = 25
$ws.Cells.Item(3,3).Value  = 25
$ws.Cells.Item(4,3).Value  = 25
$ws.Cells.Item(5,3).Value  = 25
$ws.Cells.Item(6,3).Value  = 85
$ws.Cells.Item(7,3).Value  = 85
$ws.Cells.Item(8,3).Value  = 85
$ws.Cells.Item(9,3).Value  = 85
$ws.Cells.Item(10,3).Value = 125
$ws.Cells.Item(11,3).Value = 125
$ws.Cells.Item(12,3).Value = 125
$ws.Cells.Item(13,3).Value = 125

# ---- Highlight the leakage readings taken at 25C (rows 2-5) with yellow fill ----
$ws.Range("D2:D5").Interior.Color = 65535

# ---- Page setup: portrait orientation ----
$ws.PageSetup.Orientation = 1

# ---- Selection matches the author's saved view state ----
$ws.Range("D5").Select()
